# "rebuild slides page and clearing out daily announcements"
#
# The deck is trimmed down from the full "daily announcements" deck to just
# the title slide and a single, blanked-out announcements slide:
#   - slides 3..12 (old slide IDs 287-296 / slide3.xml..slide12.xml) are
#     removed entirely.
#   - the remaining second slide ("Wednesday, January 15" with a long list
#     of announcements) is rebuilt as a fresh "Thursday, August 27" slide
#     whose body is just a placeholder line.

$p = $ppt.ActivePresentation

# Delete slides 3 through 12 (walk backwards so indices of the ones we still
# need to delete don't shift while we work).
for ($i = $p.Slides.Count; $i -ge 3; $i--) {
    $p.Slides.Item($i).Delete()
}

# Rebuild the remaining announcements slide (slide 2).
$s = $p.Slides.Item(2)

# Title placeholder.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Thursday, August 27"

# Content placeholder - replace the whole bulleted list with a single
# placeholder line.
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = [char]0x2026 + "Announcements will go here."
$body.Font.Size = 26
